$d = $word.ActiveDocument

# 1) The two runs that made up the "Under testningen ..." / "visar vikten ..."
#    sentence are replaced by a single new sentence (one run, matching the
#    target: Find&Replace naturally collapses the paragraph to one run).
$range = $d.Content
$range.Find.Execute(
    "Under testningen blev det också tydligt för oss att alla våra antaganden inte stämde, vilket visar vikten av att kontinuerligt testa våra antaganden under projektets gång.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Efterhand som vi testat våra antaganden har det blivit tydligt för oss att vi måste vara bestämda med vår kund och visa att vi inte kommer göra mer än vad som faktiskt stod i projektbeskrivningen. Vi måste tänka på att vi har andra kurser också och inte låta diverse handledare spåna på med nya idéer och krav. ",
    2
) | Out-Null

# 2) Insert three new paragraphs right before the final (trailing empty)
#    paragraph of the document body, i.e. directly after the paragraph we
#    just edited. Each inherits the same paragraph/run formatting as its
#    neighbours (ind left=360, Lato 24pt) because InsertParagraphBefore
#    clones the formatting of the paragraph it splits.
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$insertPos = $lastPara.Range.Start

$r = $d.Range($insertPos, $insertPos)
$r.InsertParagraphBefore() | Out-Null
$r.InsertParagraphBefore() | Out-Null
$r.InsertParagraphBefore() | Out-Null

# The three freshly-created empty paragraphs now sit immediately before the
# trailing empty paragraph; fill them in with the reflection text.
$paras2 = $d.Paragraphs
$n = $paras2.Count
$newPara1 = $paras2.Item($n - 3)
$newPara2 = $paras2.Item($n - 2)
$newPara3 = $paras2.Item($n - 1)

$newPara1.Range.InsertAfter("Den kommunikationen har dock gått väl och vi i gruppen känner att detta blivit mer tydligt efterhand som vi förklarat vår situation för företaget. ")
$newPara2.Range.InsertAfter("Däremot siktar vi såklart på ett väl genomfört projekt som uppfyller samtliga krav vi kommit överens om.")
$newPara3.Range.InsertAfter("Att vissa av våra antaganden inte stämde visar på vikten av att kontinuerligt testa antaganden vi gör under projektets gång.")
